$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.062.85"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "'2.935.01"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'373.15"
$ws.Range("E5").Value = "  -1.39%  "

$ws.Range("D6").Value = "'100.56"
$ws.Range("E6").Value = "  -3.91%  "

$ws.Range("D7").Value = "'0.533"
$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  -2.12%  "

$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  -2.96%  "

$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "'3.396.40"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").Value = "'17.96"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").Value = "'7.47"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").Value = "'11.24"
$ws.Range("E16").Value = "  +51.60%  "

$ws.Range("D17").Value = "'2.932.56"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").Value = "'0.972"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "'51.003.03"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").Value = "'3.14"
$ws.Range("E20").Value = "  -5.69%  "

$ws.Range("D21").Value = "'12.35"
$ws.Range("E21").Value = "  -4.36%  "

$ws.Range("D22").Value = "'0.0₃0953"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'68.58"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'264.03"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  +9.65%  "

$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").Value = "'7.35"
$ws.Range("E27").Value = "  -3.42%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "'0.111"
$ws.Range("E29").Value = "  -3.58%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'25.53"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.163"
$ws.Range("E31").Value = "  -4.20%  "

$ws.Range("D32").Value = "'9.91"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").Value = "'50.62"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'32.96"
$ws.Range("E35").Value = "  -6.04%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0440"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").Value = "'3.17"
$ws.Range("E38").Value = "  +4.76%  "

$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").Value = "'16.37"
$ws.Range("E40").Value = "  -4.64%  "

$ws.Range("E41").Value = "  -3.67%  "

$ws.Range("D42").Value = "'2.46"
$ws.Range("E42").Value = "  -4.43%  "

$ws.Range("D43").Value = "'119.72"
$ws.Range("E43").Value = "  -4.50%  "

$ws.Range("D44").Value = "'21.04"
$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("D45").Value = "'0.277"
$ws.Range("E45").Value = "  -2.96%  "

$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("D47").Value = "'3.30"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").Value = "'2.29"
$ws.Range("E48").Value = "  -3.27%  "

$ws.Range("D49").Value = "'1.980.73"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("D51").Value = "'1.30"
$ws.Range("E51").Value = "  +0.74%  "
